$wb = $excel.ActiveWorkbook

# Sheet4 = "CSC-CSCSoCECBiaSY" (index 4, onshore wind es calibration row)
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

# Row 7 = "onshore wind es" Share of existing capacity, columns B:AE -> 0.25
$ws.Range("B7:AE7").Value = 0.25
